$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.188.15'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.630.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.57%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.547'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.630.17'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.134'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.92%  '
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.348'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("E15").Value = '  +3.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.111.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.929.59'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.630.43'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '374.37'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.55%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.82%  '
$ws.Range("E28").Value = '  +2.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.760.60'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '580.35'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.41'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.86'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.52'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '158.96'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.19'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.68%  '
$ws.Range("E40").Value = '  +4.71%  '
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.36'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.12'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0₆0318'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +10.88%  '
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '155.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.71'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("E50").Value = '  -2.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.82'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.49%  '
